# Credit purchase payment upload: point the upload row at the new
# analyst's local file path and leave the selection where she left it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "C:\Users\stellah.ireri\git\MKOPA_Regression_Test_Channel\FilesToUpload\Mpesalatest File.csv"

$ws.Range("A10").Select()
